$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 3 becomes the "Sample" example: Active -> Yes, Fname -> "Upendra Kishore"
$ws.Range("A3").Value = "Yes"
$ws.Range("B3").Value = "Upendra Kishore"

# Row 4 flips to Active -> No
$ws.Range("A4").Value = "No"
